# cboc_signin_sheet.xlsx — "test of adding cboc/core row"
#
# 1. Bump the header's Month/Year label from May 2021 to June 2021.
# 2. Add a new (currently empty) row below the header, formatted with the
#    same "core row" look: bold 10pt Times New Roman text inside a thick
#    navy box border — a template row ready for sign-in entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update header text -------------------------------------------------
$ws.Range("A1").Value = "Month/Year: June 2021"

# --- 2. New formatted core row ----------------------------------------------
$coreRow = $ws.Range("A2")

$coreRow.Font.Name = "Times New Roman"
$coreRow.Font.Bold = $true
$coreRow.Font.Size = 10

$coreRow.Borders.LineStyle = 1    # xlContinuous
$coreRow.Borders.Weight = 4       # xlThick
$coreRow.Borders.Color = 5512192  # BGR for RGB(00,1C,54) navy, matches header box border

Write-Output "Updated Month/Year header and added formatted core row at A2"
